$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 (Jill): LastName grant -> Grant ; EmailAddress -> jill@give.com
$ws.Range("D3").Value = "Grant"
$ws.Range("F3").Value = "jill@give.com"

# Row 4: FirstName Bob -> John ; EmailAddress -> john@give.com ; PostalCode -> W1 3QP
$ws.Range("C4").Value = "John"
$ws.Range("E4").Value = "John"
$ws.Range("F4").Value = "john@give.com"
$ws.Range("K4").Value = "W1 3QP"

# Row 5: Title Mr -> Miss ; FirstName Bob -> Cathy ; LastName Smith -> Holmes ; EmailAddress -> cathy@give.com ; PostalCode -> W1 4QP
$ws.Range("B5").Value = "Miss"
$ws.Range("C5").Value = "Cathy"
$ws.Range("D5").Value = "Holmes"
$ws.Range("E5").Value = "Cathy"
$ws.Range("F5").Value = "cathy@give.com"
$ws.Range("K5").Value = "W1 4QP"

# Row 6: FirstName Bob -> Brian ; LastName Smith -> Monroe ; EmailAddress -> brian@give.com ; PostalCode -> W1 5QP
$ws.Range("C6").Value = "Brian"
$ws.Range("D6").Value = "Monroe"
$ws.Range("E6").Value = "Brian"
$ws.Range("F6").Value = "brian@give.com"
$ws.Range("K6").Value = "W1 5QP"

# Update selection to match the new view
$ws.Range("B2:K6").Select()
